# Swap the presentation's theme colour scheme (ppt/theme/theme1.xml) from the
# "Integral" / "Red Violet" palette to the default Office "Office Theme" palette.
#
# PowerPoint stores RGB() values as a BGR-packed long (0xBBGGRR), so convert
# each "RRGGBB" hex string accordingly before assigning ThemeColor.RGB.
function Hex2BGR($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# Standard ThemeColorScheme slot order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink (indices 1-12).
$officeTheme = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $tcs.Item($i).RGB = Hex2BGR $officeTheme[$i - 1]
}
